$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 with the same header style as A1:C1
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Update the computed metric values in B2 and C2
$ws.Range("B2").Value = 0.05579388202094215
$ws.Range("C2").Value = 0.9992445852779073

# Add the new "Tipo" data cell D2 (plain style, matching A2/B2/C2)
$ws.Range("D2").Value = "single"
